$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Genome / Guide Sequence header rows for the four new blocks ---
$ws.Range("A19").Value2 = "Genome"
$ws.Range("B19").Value2 = "Frock_EMX"
$ws.Range("A20").Value2 = "Guide Sequence"
$ws.Range("B20").Value2 = "GAGTCCGAGCAGAAGAAGAAGGG"

$ws.Range("A37").Value2 = "Genome"
$ws.Range("B37").Value2 = "Frock_RAG1A"
$ws.Range("A38").Value2 = "Guide Sequence"
$ws.Range("B38").Value2 = "GCCTCTTTCCCACCCACCTTGGG"

# --- Frock_EMX similarity rows ---
$ws.Range("B21").Value2 = "GTGGCGCATTGCCACGAAGCAGG"
$ws.Range("C21").Value2 = -0.32079273748575998
$ws.Range("B22").Value2 = "CCAGCGACGTGCCCCAGGACGGG"
$ws.Range("C22").Value2 = -0.32079273748575998
$ws.Range("B23").Value2 = "TGCGGAGGGGAGTGGACTTAGGG"
$ws.Range("C23").Value2 = -0.30839006191383
$ws.Range("B24").Value2 = "CGTGGGCCCAAGCTGGACTCTGG"
$ws.Range("C24").Value2 = -0.33592883798290502

# --- Frock_RAG1A similarity rows ---
$ws.Range("B39").Value2 = "CACATATTAAATTTTCAGAATGG"
$ws.Range("C39").Value2 = -0.38215702939849999
$ws.Range("B39").Font.Name = "Arial Unicode MS"
$ws.Range("B39").Font.Size = 10
$ws.Range("B39").VerticalAlignment = -4108
$ws.Range("B40").Value2 = "TCAGGCAAGGATCAGCAGCAAGG"
$ws.Range("C40").Value2 = -0.37908082937397403
$ws.Range("B41").Value2 = "CTCAGATGCCTCAAAGTCATGGG"
$ws.Range("C41").Value2 = -0.36456938458516402
$ws.Range("B42").Value2 = "CTTGTTCCTGCTGGCTCTGAGGG"
$ws.Range("C42").Value2 = -0.35465526773894401

# --- Wang_WAS-CR5 block (header + similarity rows) ---
$ws.Range("A28").Value2 = "Genome"
$ws.Range("B28").Value2 = "Wang_WAS-CR5"
$ws.Range("A29").Value2 = "Guide Sequence"
$ws.Range("B29").Value2 = "CCCTGTGTCTCTGGATGGATGGG"
$ws.Range("B30").Value2 = "ACTAGTGAATGAAACTGCAGAGG"
$ws.Range("C30").Value2 = -0.40116985011574902
$ws.Range("B31").Value2 = "TTTGTGCTTATCTTAATACCAGG"
$ws.Range("C31").Value2 = -0.40331181544619299
$ws.Range("B32").Value2 = "AACCCCCCCAGGTTACCTGTGGG"
$ws.Range("C32").Value2 = -0.403228086189868
$ws.Range("B33").Value2 = "CGGGTGGATCACCTGAGGTCAGG"
$ws.Range("C33").Value2 = -0.38898218306267601

# --- Highlight the best Cho_ccr-1 match and note it ---
$ws.Range("C11").Value2 = "more than least similarity"
$ws.Range("B11").Interior.Color = 14123647

# --- View / print setup ---
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("C8").Select()
$ws.PageSetup.Orientation = 1
